# Weekly data update: a new price record is inserted as row 127
# (pushing all existing rows 127:224 down by one, to 128:225).
# This matches the canonical OOXML diff, where the row that used to be at
# position N now lives at N+1 for every N in 127..224, and a brand new
# record appears at row 127.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 127; Excel shifts rows 127:224 down to 128:225,
# carrying all of their existing values/formatting with them.
$ws.Rows("127:127").Insert()

# Populate the newly inserted row 127 with the new weekly record.
$ws.Range("A127").Value = 11
$ws.Range("B127").Value = "Vega Monumental Concepción"
$ws.Range("C127").Value = "Bíobío"
$ws.Range("D127").Value = 45062
$ws.Range("E127").Value = 8
$ws.Range("F127").Value = 100112032
$ws.Range("G127").Value = "Zapallo italiano"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 100
$ws.Range("K127").Value = 12000
$ws.Range("L127").Value = 13000
$ws.Range("M127").Value = 12500
$ws.Range("N127").Value = "`$/caja 50 unidades"
$ws.Range("O127").Value = "Región de Arica y Parinacota"
$ws.Range("P127").Value = 250
$ws.Range("Q127").Value = 50
$ws.Range("R127").Value = "Hortaliza"
